$d = $word.ActiveDocument

# --- Change 1: merge "L" + "G 32LK6100 ..." into a single run "LG 32LK6100 ..." ---
$d.Content.Find.Execute("LG 32LK6100", $true, $false, $false, $false, $false,
                         $true, 1, $false, "LG 32LK6100", 2)

# --- Change 2 & 3: fix "HeadPhoes" -> "HeadPhones", splitting into "HeadPho" / "n" / "es"
#     runs with the (uniquely-named) _GoBack bookmark relocated between "n" and "es".
#     Because _GoBack is unique document-wide, re-adding it here removes the old
#     occurrence that used to sit after " £579" in the HP laptop paragraph. ---
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "HeadPhoes") {
        $target = $p
    }
}

$start = $target.Range.Start

# Insert the missing "n" right after "HeadPho" (7 characters in).
$insertionPoint = $d.Range($start + 7, $start + 7)
$insertionPoint.InsertAfter("n")
# The paragraph text is now "HeadPhones" as a single run.

# Force a run boundary between "HeadPho" and "n" using a transient bookmark -
# the split persists even once the bookmark itself is removed.
$splitRange = $d.Range($start + 7, $start + 7)
$d.Bookmarks.Add("TempSplit", $splitRange)
$d.Bookmarks.Item("TempSplit").Delete()

# Place the real _GoBack bookmark between "n" and "es" - this both forces the
# second run boundary and (being unique) relocates the bookmark from its old spot.
$goBackRange = $d.Range($start + 8, $start + 8)
$d.Bookmarks.Add("_GoBack", $goBackRange)
